$wb = $excel.ActiveWorkbook

# --- P2G sheet: convert CAPEX, Fixed O&M and Variable O&M from absolute
#     units down to millions (divide by 1,000,000) ---
$wsP2G = $wb.Worksheets.Item("P2G")

$capexRange = $wsP2G.Range("B2:B33")
foreach ($cell in $capexRange.Cells) {
    $cell.Value = $cell.Value2 / 1000000
}
$capexRange.NumberFormat = "0.0"

$fixedOMRange = $wsP2G.Range("C2:C33")
foreach ($cell in $fixedOMRange.Cells) {
    $cell.Value = $cell.Value2 / 1000000
}

$varOMRange = $wsP2G.Range("D2:D33")
foreach ($cell in $varOMRange.Cells) {
    $cell.Value = $cell.Value2 / 1000000
}

# --- FC sheet: convert CAPEX and Fixed O&M down to millions as well ---
$wsFC = $wb.Worksheets.Item("FC")

$fcCapexRange = $wsFC.Range("B2:B33")
foreach ($cell in $fcCapexRange.Cells) {
    $cell.Value = $cell.Value2 / 1000000
}

$fcFixedOMRange = $wsFC.Range("C2:C33")
foreach ($cell in $fcFixedOMRange.Cells) {
    $cell.Value = $cell.Value2 / 1000000
}

# --- Restore the active selections on each sheet ---
$wsFC.Activate()
$wsFC.Range("G4").Select()

$wsP2G.Activate()
$wsP2G.Range("L15").Select()
